# "error solve ifrs list"
#
# The sheet holds per-year IFRS financial figures for Sewoo Global.
# Rows 2-6 (fiscal years 2014-2018) had their data columns (D:AJ)
# populated with wrong (far too large) figures, and some columns that
# should not exist for this company (column J, column O, and in a
# couple of rows S/V) need to be cleared entirely. Rows 7-9 (the three
# "(E)" estimate years) were pulled in by mistake and must have all of
# their data cells (D:AJ) wiped, leaving just the A/B/C year-label
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected per-row figures. $null entries mean "this cell must not
# exist after the edit" (ClearContents), everything else is a plain
# value write.
$rows = @{
    2 = @{
        D=445;    E=25;   F=25;   G=15;  H=3;   I=3;   J=$null
        K=362;    L=59;   M=304;  N=304; O=$null
        P=118;    Q=6;    R=-57;  S=4;   T=2;   U=4;   V=4
        W=5.58;   X=0.58; Y=0.85; Z=0.74
        AA=19.28; AB=161.72; AC=11; AD=96.45; AE=1282; AF=0.82
        AG=0;     AH=0;   AI=0;   AJ=23690705
    }
    3 = @{
        D=396;    E=16;   F=16;   G=11;  H=6;   I=6;   J=$null
        K=370;    L=53;   M=317;  N=317; O=$null
        P=118;    Q=12;   R=-11;  S=-4;  T=10;  U=2;   V=$null
        W=4;      X=1.58; Y=2.01; Z=1.71
        AA=16.71; AB=167.11; AC=26; AD=85.3; AE=1337; AF=1.68
        AG=0;     AH=0;   AI=0;   AJ=23690705
    }
    4 = @{
        D=424;    E=8;    F=8;    G=9;   H=4;   I=4;   J=$null
        K=392;    L=73;   M=319;  N=319; O=$null
        P=118;    Q=-3;   R=7;    S=$null; T=5; U=-8;  V=$null
        W=1.85;   X=1.01; Y=1.35; Z=1.13
        AA=22.98; AB=170.67; AC=18; AD=122.63; AE=1346; AF=1.65
        AG=0;     AH=0;   AI=0;   AJ=23690705
    }
    5 = @{
        D=364;    E=17;   F=17;   G=3;   H=4;   I=4;   J=$null
        K=403;    L=81;   M=322;  N=322; O=$null
        P=118;    Q=-14;  R=-29;  S=70;  T=3;   U=-17; V=67
        W=4.6;    X=1.23; Y=1.39; Z=1.12
        AA=25.13; AB=177.42; AC=19; AD=95.48; AE=1361; AF=1.32
        AG=0;     AH=0;   AI=0;   AJ=23690705
    }
    6 = @{
        D=322;    E=5;    F=5;    G=1;   H=5;   I=5
        K=434;    L=107;  M=327;  N=327
        P=118;    Q=46;   R=-17;  S=0;   T=2;   U=44;  V=70
        W=1.46;   X=1.64; Y=1.63; Z=1.26
        AA=32.9;  AB=179.16; AC=22; AD=77.37; AE=1379; AF=1.25
        AG=0;     AH=0;   AI=0;   AJ=23690705
    }
}

foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$r"
        $val = $cols[$col]
        if ($null -eq $val) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $val
        }
    }
}

# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) were mistakenly filled in
# with data that doesn't belong here; wipe the whole data range, keeping
# only the A/B/C year-label cells.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

Write-Output "applied ifrs corrections to rows 2-9"
